# "Generate Report for Handoff" — refresh the localization-status report:
# status flips from "In Translation" to "Ready for handoff" and the
# associated timestamps advance to the new handoff-generation run.

$wb = $excel.ActiveWorkbook

$newStatus       = "Ready for handoff"
$overviewDate    = "2016-08-26 14:49:57"
$zhHandoffDate   = "2016-08-26 14:49:52"

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = $overviewDate

# Column widths E:F grow to fit the new, longer status text.
$overview.Columns.Item(5).ColumnWidth = 16.333333333333332
$overview.Columns.Item(6).ColumnWidth = 16.333333333333332

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = $zhHandoffDate
$zhcn.Columns.Item(3).ColumnWidth = 16.333333333333332

# --- de-de sheet ----------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = $overviewDate
$dede.Columns.Item(3).ColumnWidth = 16.333333333333332
